$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Package / User data. Row 1 (header) stays the same.
$values = @(
    @("Package", "User"),
    @("autoawq", "A"),
    @("duckdb", "B"),
    @("pyjwt", "C"),
    @("dash-extensions", "D"),
    @("io", "E"),
    @("getpadd", "F"),
    @("jwcrypto", "G")
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i][0]
    $ws.Cells.Item($row, 2).Value = $values[$i][1]
}

# Set font size for the whole sheet to 11 (new default size for this workbook)
$ws.Cells.Font.Size = 11

# Apply left alignment to the Package column data (A2:A8)
$ws.Range("A2:A8").HorizontalAlignment = -4131

# Leave the same cell selected as in the source workbook
$ws.Range("D8").Select()
